$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update rows 31-60 on Sheet1: flip sign pattern to the "1000" group
for ($r = 31; $r -le 60; $r++) {
    $ws1.Cells.Item($r, 1).Value = 1000
    $ws1.Cells.Item($r, 2).Value = 4000
    $ws1.Cells.Item($r, 3).Value = 2
    $ws1.Cells.Item($r, 4).Value = 50
    $ws1.Cells.Item($r, 5).Value = 500
    $ws1.Cells.Item($r, 6).Value = 150
    $ws1.Cells.Item($r, 7).Value = 0
}

# Delete rows 61-100 entirely (shrinks used range / dimension to A1:G60)
$ws1.Range("A61:G100").Delete(-4162) | Out-Null

# Update selection on Sheet1
$ws1.Range("I40").Select() | Out-Null
